# Auto-generated Excel COM-interop script to apply market-price refresh diffs
# to the Aegis_Profits workbook (Sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ALC!row 19 - "Unbreak My Heart" / "Roof Tile" (item 7015)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1068.2
$ws.Range("I19").Value = 982.5714
$ws.Range("J19").Value = 1143.125
$ws.Range("K19").Value = 982.5714
$ws.Range("L19").Value = 1143.125
$ws.Range("M19").Value = -807.5714
$ws.Range("N19").Value = -1493.125

# ALC!row 62 - "The Mustache Suits Him" / "Enchanted Mythrite Ink" (item 27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2681.2856
$ws.Range("I62").Value = 1996.6666
$ws.Range("J62").Value = 3194.75
$ws.Range("K62").Value = 1996.6666
$ws.Range("L62").Value = 3194.75
$ws.Range("M62").Value = -1372.6666
$ws.Range("N62").Value = -4442.75

# ALC!row 65 - "Forgery of Convenience (L)" / "Enchanted Mythrite Ink" (item 27781)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2681.2856
$ws.Range("I65").Value = 1996.6666
$ws.Range("J65").Value = 3194.75
$ws.Range("K65").Value = 9983.333000000001
$ws.Range("L65").Value = 15973.75
$ws.Range("M65").Value = -6863.333000000001
$ws.Range("N65").Value = -22213.75

# ALC!row 76 - "Warding Off Temptation" / "Enchanted Hardsilver Ink" (item 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4213.857
$ws.Range("I76").Value = 4249.5
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 4249.5
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -3934.5
$ws.Range("N76").Value = -4630

# ALC!row 79 - "The Garden of Arcane Delights (L)" / "Enchanted Hardsilver Ink" (item 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4213.857
$ws.Range("I79").Value = 4249.5
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 4249.5
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -3157.5
$ws.Range("N79").Value = -6184

# ALC!row 86 - "Filling in the Blanks" / "Enchanted Aurum Regis Ink" (item 12603)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 73406.28999999999
$ws.Range("I86").Value = 78922.16
$ws.Range("J86").Value = 1700
$ws.Range("K86").Value = 78922.16
$ws.Range("L86").Value = 1700
$ws.Range("M86").Value = -77799.16
$ws.Range("N86").Value = -3946

# ALC!row 89 - "Ink into Antiquity (L)" / "Enchanted Aurum Regis Ink" (item 12603)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 73406.28999999999
$ws.Range("I89").Value = 78922.16
$ws.Range("J89").Value = 1700
$ws.Range("K89").Value = 394610.8
$ws.Range("L89").Value = 8500
$ws.Range("M89").Value = -388994.8
$ws.Range("N89").Value = -19732

# ALC!row 92 - "Whinier than the Sword" / "Enchanted Koppranickel Ink" (item 19901)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1041.7222
$ws.Range("I92").Value = 862.6667
$ws.Range("J92").Value = 1399.8334
$ws.Range("K92").Value = 862.6667
$ws.Range("L92").Value = 1399.8334
$ws.Range("M92").Value = 385.3333
$ws.Range("N92").Value = -3895.8334

# ALC!row 106 - "Making Your Mark" / "Enchanted Palladium Ink" (item 19903)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3000
$ws.Range("I106").Value = 3166.6667
$ws.Range("J106").Value = 2500
$ws.Range("K106").Value = 3166.6667
$ws.Range("L106").Value = 2500
$ws.Range("M106").Value = -2535.6667
$ws.Range("N106").Value = -3762

# ALC!row 116 - "Growing Up" / "Growth Formula Kappa" (item 27778)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3674.3333
$ws.Range("J116").Value = 4511.5
$ws.Range("L116").Value = 4511.5
$ws.Range("N116").Value = -11395.5

# ALC!row 129 - "Practical Command" / "Commanding Craftsman's Draught" (item 36115)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 942.9355
$ws.Range("I129").Value = 517
$ws.Range("J129").Value = 1006.03705
$ws.Range("K129").Value = 1551
$ws.Range("L129").Value = 3018.11115
$ws.Range("M129").Value = 3449
$ws.Range("N129").Value = -13018.11115

# ARM!row 28 - "246 Kinds of Cheese" / "Iron Frypan" (item 19534)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 18401.5
$ws.Range("I28").Value = 15868.667
$ws.Range("J28").Value = 26000
$ws.Range("K28").Value = 15868.667
$ws.Range("L28").Value = 26000
$ws.Range("M28").Value = -15676.667
$ws.Range("N28").Value = -26384

# ARM!row 97 - "Ore for Me" / "High Steel Ingot" (item 19941)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 48948.145
$ws.Range("I97").Value = 72279.21000000001
$ws.Range("J97").Value = 2286
$ws.Range("K97").Value = 72279.21000000001
$ws.Range("L97").Value = 2286
$ws.Range("M97").Value = -71783.21000000001
$ws.Range("N97").Value = -3278

# ARM!row 99 - "Home Cooking" / "Doman Iron Frypan" (item 19534)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 18401.5
$ws.Range("I99").Value = 15868.667
$ws.Range("J99").Value = 26000
$ws.Range("K99").Value = 15868.667
$ws.Range("L99").Value = 26000
$ws.Range("M99").Value = -12873.667
$ws.Range("N99").Value = -31990

# BSM!row 86 - "Through Thick and Thin" / "Adamantite Nugget" (item 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 52673.863
$ws.Range("I86").Value = 76364.60000000001
$ws.Range("J86").Value = 1908
$ws.Range("K86").Value = 76364.60000000001
$ws.Range("L86").Value = 1908
$ws.Range("M86").Value = -75241.60000000001
$ws.Range("N86").Value = -4154

# BSM!row 89 - "Piercing Eyes Deserve Piercing Shafts (L)" / "Adamantite Nugget" (item 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 52673.863
$ws.Range("I89").Value = 76364.60000000001
$ws.Range("J89").Value = 1908
$ws.Range("K89").Value = 381823
$ws.Range("L89").Value = 9540
$ws.Range("M89").Value = -376207
$ws.Range("N89").Value = -20772

# BSM!row 94 - "High Steal" / "High Steel Nugget" (item 19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 551.7727
$ws.Range("J94").Value = 557.2
$ws.Range("L94").Value = 557.2
$ws.Range("N94").Value = -1459.2

# CRP!row 133 - "Yimepi's Country Charms" / "Ginseng Earrings" (item 43328)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 63400
$ws.Range("J133").Value = 63400
$ws.Range("L133").Value = 63400
$ws.Range("N133").Value = -68460

# CRP!row 135 - "The Wing's Wings" / "Ceiba Wings" (item 42008)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 50861.332
$ws.Range("J135").Value = 50861.332
$ws.Range("L135").Value = 50861.332
$ws.Range("N135").Value = -61001.332

# CRP!row 138 - "Bow Out" / "Acacia Longbow" (item 42302)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 71232.22
$ws.Range("J138").Value = 71232.22
$ws.Range("L138").Value = 71232.22
$ws.Range("N138").Value = -81512.22

# CUL!row 34 - "Fever Pitch" / "Chamomile Tea" (item 4749)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 999.25
$ws.Range("J34").Value = 1249
$ws.Range("L34").Value = 3747
$ws.Range("N34").Value = -3915

# CUL!row 70 - "Persona non Gratin" / "Dhalmel Gratin" (item 12867)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 68927.60000000001
$ws.Range("I70").Value = 92901.27
$ws.Range("K70").Value = 278703.81
$ws.Range("M70").Value = -278388.81

# CUL!row 73 - "Recipe for Disaster (L)" / "Dhalmel Gratin" (item 12867)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 68927.60000000001
$ws.Range("I73").Value = 92901.27
$ws.Range("K73").Value = 278703.81
$ws.Range("M73").Value = -277611.81

# CUL!row 86 - "Let's Not Get Sappy" / "Birch Syrup" (item 12892)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 595
$ws.Range("J86").Value = 618
$ws.Range("L86").Value = 1854
$ws.Range("N86").Value = -4226

# CUL!row 89 - "Luxury Spillover (L)" / "Birch Syrup" (item 12892)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 595
$ws.Range("J89").Value = 618
$ws.Range("L89").Value = 5562
$ws.Range("N89").Value = -17418

# CUL!row 107 - "Slippery Service" / "Frantoio Oil" (item 27838)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 463864.88
$ws.Range("I107").Value = 589.0909
$ws.Range("J107").Value = 732077.2
$ws.Range("K107").Value = 1767.2727
$ws.Range("L107").Value = 2196231.6
$ws.Range("M107").Value = 152.7273
$ws.Range("N107").Value = -2200071.6

# GSM!row 18 - "Gorgeous Gorget" / "Brass Gorget" (item 4309)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# GSM!row 97 - "If I'd a Koppranickel for Every Time..." / "Koppranickel Ingot" (item 19940)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 90911300
$ws.Range("I97").Value = 125002616
$ws.Range("J97").Value = 1100
$ws.Range("K97").Value = 125002616
$ws.Range("L97").Value = 1100
$ws.Range("M97").Value = -125002120
$ws.Range("N97").Value = -2092

# LTW!row 93 - "Hide to Go Seek" / "Gagana Leather" (item 19993)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2055
$ws.Range("I93").Value = 2481.3635
$ws.Range("J93").Value = 1533.8889
$ws.Range("K93").Value = 2481.3635
$ws.Range("L93").Value = 1533.8889
$ws.Range("M93").Value = -1233.3635
$ws.Range("N93").Value = -4029.8889

# LTW!row 100 - "Tiger in the Sack" / "Tiger Leather" (item 19995)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2646.6667
$ws.Range("I100").Value = 2133.3333
$ws.Range("J100").Value = 3160
$ws.Range("K100").Value = 2133.3333
$ws.Range("L100").Value = 3160
$ws.Range("M100").Value = -1592.3333
$ws.Range("N100").Value = -4242

# WVR!row 122 - "Heavy Armoire" / "Dark Hempen Cloth" (item 36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2293.389
$ws.Range("J122").Value = 3110.889
$ws.Range("L122").Value = 9332.667000000001
$ws.Range("N122").Value = -14232.667
